# Update countries & provincias Spain
# - Republica Dominicana overtakes Panama and Paises Bajos in total cases,
#   so it moves up two positions in the (descending, sorted-by-total-cases)
#   country table. Panama's and Paises Bajos' old figures cascade down one
#   row each, and Republica Dominicana gets fresh figures.
# - Mozambique overtakes Zimbabue, so it moves up one position; Zimbabue's
#   old figures cascade down one row, Mozambique gets fresh figures.
# - A handful of other countries (Estados Unidos, Haiti, Tunez, Benin,
#   Crucero) just get refreshed totals without changing rank.
# - The "last updated" timestamp banner is bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name column (A) swaps caused by the re-sort ---
$ws.Range("A41").Value = "Republica Dominicana"
$ws.Range("A42").Value = "Panama"
$ws.Range("A43").Value = "Paises Bajos"

$ws.Range("A134").Value = "Mozambique"
$ws.Range("A135").Value = "Zimbabue"

# --- Row 4: Estados Unidos (fresh totals, no rank change) ---
$ws.Range("B4").Value = 3836674
$ws.Range("C4").Value = 3403
$ws.Range("D4").Value = 1775491
$ws.Range("E4").Value = 1918290
$ws.Range("G4").Value = 16
$ws.Range("H4").Value = 142893

# --- Row 41: now Republica Dominicana (fresh totals) ---
$ws.Range("B41").Value = 52855
$ws.Range("C41").Value = 1336
$ws.Range("D41").Value = 25094
$ws.Range("E41").Value = 26780
$ws.Range("G41").Value = 10
$ws.Range("H41").Value = 981

# --- Row 42: now Panama (old Republica Dominicana... no: cascaded old row41 data) ---
$ws.Range("B42").Value = 52261
$ws.Range("C42").Value = 0
$ws.Range("D42").Value = 27494
$ws.Range("E42").Value = 23696
$ws.Range("H42").Value = 1071

# --- Row 43: now Paises Bajos (cascaded old row42 data) ---
$ws.Range("B43").Value = 51725
$ws.Range("C43").Value = 144
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("H43").Value = 6136

# --- Row 91: Haiti (fresh totals, no rank change) ---
$ws.Range("B91").Value = 6878
$ws.Range("C91").Value = 44
$ws.Range("D91").Value = 5578
$ws.Range("E91").Value = 1243

# --- Row 134: now Mozambique (fresh totals) ---
$ws.Range("B134").Value = 1491
$ws.Range("C134").Value = 56
$ws.Range("D134").Value = 472
$ws.Range("E134").Value = 1009
$ws.Range("H134").Value = 10

# --- Row 135: now Zimbabue (cascaded old row134 data) ---
$ws.Range("B135").Value = 1478
$ws.Range("C135").Value = 0
$ws.Range("D135").Value = 439
$ws.Range("E135").Value = 1014
$ws.Range("H135").Value = 25

# --- Row 136: Tunez (fresh totals, no rank change) ---
$ws.Range("B136").Value = 1374
$ws.Range("C136").Value = 26
$ws.Range("D136").Value = 1097
$ws.Range("E136").Value = 227

# --- Row 157: Crucero (fresh totals, no rank change) ---
$ws.Range("B157").Value = 628
$ws.Range("C157").Value = 4
$ws.Range("E157").Value = 153

# --- Banner timestamp bump ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Julio de 2020 a las 16:58"
